# Score_iterations.xlsx - add new GFG log entries (stacks topic)
# "return max possible sum equal in 3 stacks"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 142: Find if an expression has duplicate parenthesis or not
$ws.Range("B142").Value = "GFG"
$ws.Range("C142").Value = "Find if an expression has duplicate parenthesis or not"
$ws.Range("D142").Value = "not done"
$ws.Range("E142").Value = "12:20 - 12:50"
$ws.Hyperlinks.Add($ws.Range("F142"), "https://www.geeksforgeeks.org/find-expression-duplicate-parenthesis-not/")
$ws.Range("F142").Style = "Hyperlink"

# Row 143: Find maximum sum possible equal sum of three stacks
$ws.Range("B143").Value = "GFG"
$ws.Range("C143").Value = "Find maximum sum possible equal sum of three stacks"
$ws.Range("D143").Value = "done 1 logical comparision error"
$ws.Range("E143").Value = "1:07 - 1:29"
$ws.Hyperlinks.Add($ws.Range("F143"), "https://www.geeksforgeeks.org/find-maximum-sum-possible-equal-sum-three-stacks/")
$ws.Range("F143").Style = "Hyperlink"

# Row 144: follow-up note row (no topic/link, just status + time)
$ws.Range("D144").Value = "coding"
$ws.Range("E144").Value = "x - 1:58"

# Row 145: Sort a stack using a temporary stack
$ws.Range("B145").Value = "GFG"
$ws.Range("C145").Value = "Sort a stack using a temporary stack"
$ws.Range("E145").Value = "12:07 - 12:30"
$ws.Hyperlinks.Add($ws.Range("F145"), "https://www.geeksforgeeks.org/sort-stack-using-temporary-stack/")
$ws.Range("F145").Style = "Hyperlink"

# Move the on-screen selection/scroll position to follow the newly typed rows
$excel.ActiveWindow.ScrollRow = 135
$ws.Range("E146").Select()
